$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODAY")

$ws.Range("K3").Value = 19
$ws.Range("L3").Value = 31.02
$ws.Range("K5").Value = 80
$ws.Range("L5").Value = 194.87
